# chore: update epidemiological weeks to include 2024
#
# The "Região" column (A) was only populated from row 54 onward (the 2025
# weekly-data block). This fills in the same "Sudeste" region label for the
# 2024 weekly-data block in rows 2-53, so the 2024 rows carry the region
# just like every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 53 (the 2024 epidemiological-week block) get column A
# filled in with the region name, matching the value already used by every
# other populated row in column A ("Sudeste").
$ws.Range("A2:A53").Value = "Sudeste"

# Restore the on-screen selection to where the author was working.
$ws.Range("I9").Select() | Out-Null

# Scroll the window so row 11 is at the top of the viewport (best effort —
# mirrors the author's recorded scroll position).
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
